$wb = $excel.ActiveWorkbook

# --- Sheet "variables" ---
$ws2 = $wb.Worksheets.Item("variables")

# Insert a "slope" row before the mean_temperature row (current row 3)
$ws2.Rows.Item(3).Insert()
$ws2.Range("A3").Value = "slope"
$ws2.Range("B3").Value = "stand"
$ws2.Range("C3").Value = "numeric"
$ws2.Range("D3").Value = "º"
$ws2.Range("E3").Value = "slope in degrees"

# Insert a "mean_precipitation" row before the canopy_cover row (now row 5)
$ws2.Rows.Item(5).Insert()
$ws2.Range("A5").Value = "mean_precipitation"
$ws2.Range("B5").Value = "stand"
$ws2.Range("C5").Value = "numeric"
$ws2.Range("D5").Value = "mm"
$ws2.Range("E5").Value = "mean precipitation for the time period "

$ws2.Range("E6").Select()

# --- Sheet "indicators" ---
$ws1 = $wb.Worksheets.Item("indicators")

$ws1.Range("D2").Value = "mean_precipitation, canopy_cover"
$ws1.Range("H4").Value = "state, dbh, n"
$ws1.Range("H3").Value = "state, dbh, h, n"

$ws1.Range("H4").Select()
